$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 3.25
$ws.Range("I4").Value = 4.1
$ws.Range("J4").Value = 2.75
$ws.Range("O4").Value = 1.5
$ws.Range("P4").Value = 2.5
$ws.Range("Q4").Value = 2.5
$ws.Range("R4").Value = 1.5
$ws.Range("S4").Value = 1.57
$ws.Range("T4").Value = 2.25
$ws.Range("X4").Value = 8
$ws.Range("Z4").Value = 17
$ws.Range("AG4").Value = 8.5
$ws.Range("AH4").Value = 19
$ws.Range("AJ4").Value = 41
$ws.Range("AO4").Value = 12
$ws.Range("AT4").Value = 2.25
$ws.Range("BD4").Value = 126

# Row 6
$ws.Range("M6").Value = 1.07
$ws.Range("N6").Value = 9
$ws.Range("Q6").Value = 2.25
$ws.Range("R6").Value = 1.62
$ws.Range("Z6").Value = 21
$ws.Range("AA6").Value = 21
$ws.Range("AG6").Value = 8
$ws.Range("AJ6").Value = 34

# Row 7
$ws.Range("G7").Value = 3.3
$ws.Range("H7").Value = 3.6
$ws.Range("I7").Value = 2.1
$ws.Range("K7").Value = 2.3
$ws.Range("U7").Value = 1.57
$ws.Range("V7").Value = 2.25
$ws.Range("W7").Value = 13
$ws.Range("X7").Value = 19
$ws.Range("AK7").Value = 15
$ws.Range("AL7").Value = 21
$ws.Range("AM7").Value = 126
$ws.Range("AS7").Value = 126
$ws.Range("AZ7").Value = 34
$ws.Range("BB7").Value = 101

# Row 8
$ws.Range("H8").Value = 3.3
$ws.Range("I8").Value = 2
$ws.Range("K8").Value = 1.95
$ws.Range("M8").Value = 1.1
$ws.Range("N8").Value = 7
$ws.Range("O8").Value = 1.5
$ws.Range("P8").Value = 2.63
$ws.Range("Q8").Value = 2.45
$ws.Range("R8").Value = 1.46
$ws.Range("S8").Value = 1.53
$ws.Range("T8").Value = 2.38
$ws.Range("U8").Value = 2.1
$ws.Range("V8").Value = 1.67
$ws.Range("W8").Value = 8.5
$ws.Range("Y8").Value = 15
$ws.Range("AA8").Value = 41
$ws.Range("AC8").Value = 7
$ws.Range("AF8").Value = 81
$ws.Range("AG8").Value = 5.5
$ws.Range("AH8").Value = 8
$ws.Range("AI8").Value = 9.5
$ws.Range("AL8").Value = 41
$ws.Range("AP8").Value = 41
$ws.Range("AS8").Value = 351
$ws.Range("AT8").Value = 2.38
$ws.Range("AX8").Value = 12
$ws.Range("BB8").Value = 251

# Row 9
$ws.Range("G9").Value = 1.85
$ws.Range("I9").Value = 4.33
$ws.Range("J9").Value = 2.5
$ws.Range("O9").Value = 1.3
$ws.Range("P9").Value = 3.5
$ws.Range("Q9").Value = 2.01
$ws.Range("R9").Value = 1.89
$ws.Range("AG9").Value = 12
$ws.Range("AH9").Value = 21
$ws.Range("AJ9").Value = 41
$ws.Range("AO9").Value = 10

# Row 10
$ws.Range("Q10").Value = 2.2
$ws.Range("R10").Value = 1.57

# Row 11
$ws.Range("U11").Value = 2.37
$ws.Range("V11").Value = 1.5

# Row 12
$ws.Range("V12").Value = 1.63

